# Trading Journal "variable names" workbook — revision 2
# Adds a new "Input or output" column (F) for the existing param rows,
# and appends 7 new "output" parameter rows (Margin value, RRR, Market
# value, Target dollar/percent value, Stop loss dollar/percent value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column F: "Input or output" header + "Input " for existing rows ---
$ws.Range("F2").Value = "Input or output"
$ws.Range("F3").Value = "Input "
$ws.Range("F4").Value = "Input "
$ws.Range("F5").Value = "Input "
$ws.Range("F6").Value = "Input "
$ws.Range("F7").Value = "Input "
$ws.Range("F8").Value = "Input "
$ws.Range("F9").Value = "Input "
$ws.Range("F10").Value = "Input "
$ws.Range("F11").Value = "Input "
$ws.Range("F12").Value = "Input "
$ws.Range("F13").Value = "Input "
$ws.Range("F14").Value = "Input "
$ws.Range("F15").Value = "Input "

# --- New row 16: Margin value ---
$ws.Range("B16").Value = "Margin value"
$ws.Range("C16").Value = "margin"
$ws.Range("E16").Value = "float"
$ws.Range("F16").Value = "output"

# --- New row 17: RRR ---
$ws.Range("B17").Value = "RRR"
$ws.Range("C17").Value = "Rrr"
$ws.Range("E17").Value = "float "
$ws.Range("F17").Value = "output"

# --- New row 18: Market value ---
$ws.Range("B18").Value = "Market value"
$ws.Range("C18").Value = "marketVal"
$ws.Range("E18").Value = "float"
$ws.Range("F18").Value = "output"

# --- New rows 19-22: Target / Stop loss dollar & percent values ---
$ws.Range("C19").Value = "targetVal"
$ws.Range("B19").Value = "Target dollar value"
$ws.Range("B20").Value = "Target percent value"
$ws.Range("B21").Value = "Stop loss  dollar value"
$ws.Range("B22").Value = "Stop loss percent value"

$ws.Range("E19").Value = "float "
$ws.Range("F19").Value = "output"
$ws.Range("E20").Value = "float "
$ws.Range("F20").Value = "output"
$ws.Range("E21").Value = "float "
$ws.Range("F21").Value = "output"
$ws.Range("E22").Value = "float "
$ws.Range("F22").Value = "output"

$ws.Range("C20").Value = "tpercentVal"
$ws.Range("C21").Value = "stoplossVal"
$ws.Range("C22").Value = "percentVal"

# --- Column B widens to fit the new longer labels ---
$ws.Columns.Item(2).AutoFit()

# --- Final selection moves to D22 ---
$null = $ws.Range("D22").Select()
